$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
$ws.Range("A8").Characters(21, 2).Text = "17"
$ws.Range("C9").Characters(27, 9).Text = "4/24/2023"
$ws.Range("C9").Characters(47, 9).Text = "4/30/2023"

# --- Cells flipping from numeric back to "n/a" text placeholders ---
$ws.Range("C23").Copy($ws.Range("G30"))
$ws.Range("E14").Copy($ws.Range("H30"))

# --- Cells flipping from "n/a" text placeholders to numeric values ---
$ws.Range("C16").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("C16").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("E16").Copy($ws.Range("E22"))
$ws.Range("E22").Value = 0
$ws.Range("C16").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 1
$ws.Range("C16").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("E16").Copy($ws.Range("E26"))
$ws.Range("E26").Value = 0
$ws.Range("C16").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 2
$ws.Range("E16").Copy($ws.Range("E27"))
$ws.Range("E27").Value = 0
$ws.Range("C16").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1
$ws.Range("C16").Copy($ws.Range("F28"))
$ws.Range("F28").Value = 1
$ws.Range("C16").Copy($ws.Range("I28"))
$ws.Range("I28").Value = 1
$ws.Range("C16").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 1
$ws.Range("C16").Copy($ws.Range("F29"))
$ws.Range("F29").Value = 1
$ws.Range("C16").Copy($ws.Range("I29"))
$ws.Range("I29").Value = 1

# --- Plain numeric value updates ---
$ws.Range("I15").Value = 15
$ws.Range("K15").Value = 200
$ws.Range("L15").Value = 650
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 650
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 81.818181818181
$ws.Range("I16").Value = 83
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 38.333333333333
$ws.Range("L16").Value = 151.515151515152
$ws.Range("M16").Value = 48.214285714285
$ws.Range("N16").Value = -78.215223097112
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 400
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -26.315789473684
$ws.Range("I17").Value = 66
$ws.Range("J17").Value = 66
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 17.857142857142
$ws.Range("M17").Value = 53.488372093023
$ws.Range("N17").Value = -16.455696202531
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = -66.666666666666
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -55.555555555555
$ws.Range("I18").Value = 66
$ws.Range("J18").Value = 69
$ws.Range("K18").Value = -4.347826086956
$ws.Range("M18").Value = -20.481927710843
$ws.Range("N18").Value = -85.652173913043
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -8.333333333333
$ws.Range("F19").Value = 58
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = 11.538461538461
$ws.Range("I19").Value = 245
$ws.Range("J19").Value = 223
$ws.Range("K19").Value = 9.865470852017
$ws.Range("L19").Value = 85.60606060606
$ws.Range("M19").Value = 66.666666666666
$ws.Range("N19").Value = -6.130268199233
$ws.Range("C20").Value = 11
$ws.Range("E20").Value = 175
$ws.Range("F20").Value = 29
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 123.076923076923
$ws.Range("I20").Value = 76
$ws.Range("J20").Value = 54
$ws.Range("K20").Value = 40.74074074074
$ws.Range("L20").Value = 46.153846153846
$ws.Range("M20").Value = -3.79746835443
$ws.Range("N20").Value = -87.878787878787
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 23.076923076923
$ws.Range("F21").Value = 134
$ws.Range("G21").Value = 123
$ws.Range("H21").Value = 8.943089430894
$ws.Range("I21").Value = 551
$ws.Range("J21").Value = 477
$ws.Range("K21").Value = 15.513626834381
$ws.Range("L21").Value = 72.1875
$ws.Range("M21").Value = 32.771084337349
$ws.Range("N21").Value = -69.641873278236
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 28
$ws.Range("J22").Value = 26
$ws.Range("K22").Value = 7.692307692307
$ws.Range("L22").Value = 115.384615384615
$ws.Range("M22").Value = 115.384615384615
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -25.925925925925
$ws.Range("F24").Value = 100
$ws.Range("G24").Value = 106
$ws.Range("H24").Value = -5.66037735849
$ws.Range("I24").Value = 509
$ws.Range("J24").Value = 396
$ws.Range("K24").Value = 28.535353535353
$ws.Range("L24").Value = 45.428571428571
$ws.Range("M24").Value = 73.720136518771
$ws.Range("C25").Value = 9
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 43
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = -12.244897959183
$ws.Range("I25").Value = 179
$ws.Range("J25").Value = 194
$ws.Range("K25").Value = -7.731958762886
$ws.Range("L25").Value = 40.944881889763
$ws.Range("M25").Value = 9.815950920245
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 17
$ws.Range("J26").Value = 6
$ws.Range("K26").Value = 183.333333333333
$ws.Range("L26").Value = 183.333333333333
$ws.Range("C27").Value = 2
$ws.Range("I27").Value = 30
$ws.Range("J27").Value = 26
$ws.Range("K27").Value = 15.384615384615
$ws.Range("L27").Value = 114.285714285714
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 1
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -75
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = -90.90909090909
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 1
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -66.666666666666
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -90.90909090909
$ws.Range("L30").Value = -33.333333333333
